# ContextFreeSQL_Video.docx - "CSV compare html updates"
#
# This script restructures the "reporting" section of the demo transcript:
#  - splits the old "look at reporting..." paragraph into three paragraphs
#    (one of which is new content, bordered)
#  - removes several now-obsolete paragraphs (old data-entities Q&A block)
#  - rewrites the "Another feature" paragraph and splits it into two
#    paragraphs (CSV feature description + new "separate video" blurb)
#  - appends a new "Text in video:" bullet list at the end of the document

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Split "...look at reporting. We have an output subdir created here
#    )this can be configured in the json file of course) and it has full
#    reporting on the state of the database." into three paragraphs.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    ' look at reporting. We have an output subdir created here )this can be configured in the json file of course) and it has full reporting on the state of the database.',
    $true, $false, $false, $false, $false, $true, 1, $false,
    ' look at reporting. ^pHere are various changes I made in the database again. Added a column on students, updated and deleted a row on studentsgrades, also deleted and added a row on students.^pNow we have an output subdirectory created here (this can be configured in the JSON file of course) and it has full reporting on the state of the database. When we ran the script before, these were generated. Let’s open the main file, database_report.html.',
    2) | Out-Null

# ---------------------------------------------------------------------
# 2) Remove the now-obsolete paragraphs:
#      "Lets run all our demo changes again. Then run our script."
#      "If it’s a data entities, I get this screen"
#      "RN: what data differences I have?"
#      "Now let me just state ... check it out in my channel"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    'Lets run all our demo changes again. Then run our script.^p',
    $true, $false, $false, $false, $false, $true, 1, $false,
    '', 2) | Out-Null

$d.Content.Find.Execute(
    'If it’s a data entities, I get this screen^p',
    $true, $false, $false, $false, $false, $true, 1, $false,
    '', 2) | Out-Null

$d.Content.Find.Execute(
    'RN: what data differences I have?^p',
    $true, $false, $false, $false, $false, $true, 1, $false,
    '', 2) | Out-Null

$d.Content.Find.Execute(
    'Now let me just state that this data comparison form can be used as a stand alone data comparer as well, to compare any tabular data. I got a video on that too, check it out in my channel^p',
    $true, $false, $false, $false, $false, $true, 1, $false,
    '', 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Rewrite the bordered "Another feature i: ..." paragraph with the new
#    "This data differences page is a feature-rich page..." text, and
#    split off a new "Another feature: ..." paragraph after it.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    'Another feature i: I can export CSVs of the data, if the size is too big to be within the SQL script. The script then will always compare the state of the data against these CSVs',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'This data differences page is a feature-rich page in which we can search and filter throughout the data differences report. It can also be used as a stand-alone data comparer as well, to compare any tabular data – not just in the context of context free SQL scripting. I got a separate video focusing just on this page which you find on my channel^pAnother feature: I can export CSV of the data, if the size is too big to be within the SQL script. The script then will always compare the state of the data against these CSV',
    2) | Out-Null

Write-Host "Step 3 done, paragraph count: $($d.Paragraphs.Count)"
